$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bullet = [char]0x25CF

# Copy formatting (styles only) from the last existing row (293) onto the
# new block of rows (294:312) so the new rows reuse the existing cellXfs
# entries instead of Excel minting duplicate style records.
$fmtSrc = $ws.Range("A293:M293")
$fmtDst = $ws.Range("A294:M312")
$fmtSrc.Copy()
$fmtDst.PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the row height used by every other data row in the sheet.
$ws.Range("A294:A312").EntireRow.RowHeight = 16

# Row 294
$ws.Range("A294").Value = $bullet
$ws.Range("B294").Value = ''
$ws.Range("C294").Value = ''
$ws.Range("D294").Value = '267'
$ws.Range("E294").Value = 'MIC'
$ws.Range("F294").Value = '3: 532'
$ws.Range("G294").Value = '3: 535'
$ws.Range("H294").Value = 0
$ws.Range("I294").Value = ' 128'
$ws.Range("J294").Value = 3
$ws.Range("K294").Value = 0.01181381428683941
$ws.Range("L294").Value = 'chen'
$ws.Range("M294").Value = '1/23/19 14:02:45'

# Row 295
$ws.Range("A295").Value = $bullet
$ws.Range("B295").Value = ''
$ws.Range("C295").Value = ''
$ws.Range("D295").Value = '267'
$ws.Range("E295").Value = 'MIC'
$ws.Range("F295").Value = '3: 408'
$ws.Range("G295").Value = '3: 411'
$ws.Range("H295").Value = 0
$ws.Range("I295").Value = '>512'
$ws.Range("J295").Value = 4
$ws.Range("K295").Value = 0.01575175238245255
$ws.Range("L295").Value = 'chen'
$ws.Range("M295").Value = '1/23/19 14:05:48'

# Row 296
$ws.Range("A296").Value = $bullet
$ws.Range("B296").Value = ''
$ws.Range("C296").Value = ''
$ws.Range("D296").Value = '267'
$ws.Range("E296").Value = 'MIC'
$ws.Range("F296").Value = '3: 432'
$ws.Range("G296").Value = '3: 435'
$ws.Range("H296").Value = 0
$ws.Range("I296").Value = '>512'
$ws.Range("J296").Value = 4
$ws.Range("K296").Value = 0.01575175238245255
$ws.Range("L296").Value = 'chen'
$ws.Range("M296").Value = '1/23/19 14:07:31'

# Row 297
$ws.Range("A297").Value = $bullet
$ws.Range("B297").Value = ''
$ws.Range("C297").Value = ''
$ws.Range("D297").Value = '3882'
$ws.Range("E297").Value = 'MIC'
$ws.Range("F297").Value = '3: 2627'
$ws.Range("G297").Value = '3: 2630'
$ws.Range("H297").Value = 0
$ws.Range("I297").Value = '>256'
$ws.Range("J297").Value = 4
$ws.Range("K297").Value = 0.017135023989033587
$ws.Range("L297").Value = 'chen'
$ws.Range("M297").Value = '1/31/19 10:21:49'

# Row 298
$ws.Range("A298").Value = $bullet
$ws.Range("B298").Value = ''
$ws.Range("C298").Value = ''
$ws.Range("D298").Value = '3882'
$ws.Range("E298").Value = 'MIC'
$ws.Range("F298").Value = '3: 2660'
$ws.Range("G298").Value = '3: 2663'
$ws.Range("H298").Value = 0
$ws.Range("I298").Value = '>256'
$ws.Range("J298").Value = 4
$ws.Range("K298").Value = 0.017135023989033587
$ws.Range("L298").Value = 'chen'
$ws.Range("M298").Value = '1/31/19 10:21:53'

# Row 299
$ws.Range("A299").Value = $bullet
$ws.Range("B299").Value = ''
$ws.Range("C299").Value = ''
$ws.Range("D299").Value = '3882'
$ws.Range("E299").Value = 'MIC'
$ws.Range("F299").Value = '3: 2732'
$ws.Range("G299").Value = '3: 2735'
$ws.Range("H299").Value = 0
$ws.Range("I299").Value = '>256'
$ws.Range("J299").Value = 4
$ws.Range("K299").Value = 0.017135023989033587
$ws.Range("L299").Value = 'chen'
$ws.Range("M299").Value = '1/31/19 10:21:57'

# Row 300
$ws.Range("A300").Value = $bullet
$ws.Range("B300").Value = ''
$ws.Range("C300").Value = ''
$ws.Range("D300").Value = '3882'
$ws.Range("E300").Value = 'MIC'
$ws.Range("F300").Value = '3: 2787'
$ws.Range("G300").Value = '3: 2790'
$ws.Range("H300").Value = 0
$ws.Range("I300").Value = '>256'
$ws.Range("J300").Value = 4
$ws.Range("K300").Value = 0.017135023989033587
$ws.Range("L300").Value = 'chen'
$ws.Range("M300").Value = '1/31/19 10:22:02'

# Row 301
$ws.Range("A301").Value = $bullet
$ws.Range("B301").Value = ''
$ws.Range("C301").Value = ''
$ws.Range("D301").Value = '3882'
$ws.Range("E301").Value = 'MIC'
$ws.Range("F301").Value = '3: 2822'
$ws.Range("G301").Value = '3: 2825'
$ws.Range("H301").Value = 0
$ws.Range("I301").Value = '>156'
$ws.Range("J301").Value = 4
$ws.Range("K301").Value = 0.017135023989033587
$ws.Range("L301").Value = 'chen'
$ws.Range("M301").Value = '1/31/19 10:22:15'

# Row 302
$ws.Range("A302").Value = $bullet
$ws.Range("B302").Value = ''
$ws.Range("C302").Value = ''
$ws.Range("D302").Value = '3882'
$ws.Range("E302").Value = 'MIC'
$ws.Range("F302").Value = '3: 2852'
$ws.Range("G302").Value = '3: 2855'
$ws.Range("H302").Value = 0
$ws.Range("I302").Value = '>256'
$ws.Range("J302").Value = 4
$ws.Range("K302").Value = 0.017135023989033587
$ws.Range("L302").Value = 'chen'
$ws.Range("M302").Value = '1/31/19 10:22:24'

# Row 303
$ws.Range("A303").Value = $bullet
$ws.Range("B303").Value = ''
$ws.Range("C303").Value = ''
$ws.Range("D303").Value = '3882'
$ws.Range("E303").Value = 'MIC'
$ws.Range("F303").Value = '3: 2878'
$ws.Range("G303").Value = '3: 2881'
$ws.Range("H303").Value = 0
$ws.Range("I303").Value = '>256'
$ws.Range("J303").Value = 4
$ws.Range("K303").Value = 0.017135023989033587
$ws.Range("L303").Value = 'chen'
$ws.Range("M303").Value = '1/31/19 10:22:36'

# Row 304
$ws.Range("A304").Value = $bullet
$ws.Range("B304").Value = ''
$ws.Range("C304").Value = ''
$ws.Range("D304").Value = '3882'
$ws.Range("E304").Value = 'MIC'
$ws.Range("F304").Value = '3: 2904'
$ws.Range("G304").Value = '3: 2907'
$ws.Range("H304").Value = 0
$ws.Range("I304").Value = '>256'
$ws.Range("J304").Value = 4
$ws.Range("K304").Value = 0.017135023989033587
$ws.Range("L304").Value = 'chen'
$ws.Range("M304").Value = '1/31/19 10:22:58'

# Row 305
$ws.Range("A305").Value = $bullet
$ws.Range("B305").Value = ''
$ws.Range("C305").Value = ''
$ws.Range("D305").Value = '3882'
$ws.Range("E305").Value = 'MIC'
$ws.Range("F305").Value = '3: 3038'
$ws.Range("G305").Value = '3: 3041'
$ws.Range("H305").Value = 0
$ws.Range("I305").Value = '>256'
$ws.Range("J305").Value = 4
$ws.Range("K305").Value = 0.017135023989033587
$ws.Range("L305").Value = 'chen'
$ws.Range("M305").Value = '1/31/19 10:23:04'

# Row 306
$ws.Range("A306").Value = $bullet
$ws.Range("B306").Value = ''
$ws.Range("C306").Value = ''
$ws.Range("D306").Value = '3882'
$ws.Range("E306").Value = 'Drug Resisted'
$ws.Range("F306").Value = '3: 3095'
$ws.Range("G306").Value = '3: 3096'
$ws.Range("H306").Value = 0
$ws.Range("I306").Value = '16'
$ws.Range("J306").Value = 2
$ws.Range("K306").Value = 0.008567511994516793
$ws.Range("L306").Value = 'chen'
$ws.Range("M306").Value = '1/31/19 10:23:09'

# Row 307
$ws.Range("A307").Value = $bullet
$ws.Range("B307").Value = ''
$ws.Range("C307").Value = ''
$ws.Range("D307").Value = '4519'
$ws.Range("E307").Value = 'MIC'
$ws.Range("F307").Value = '3: 3605'
$ws.Range("G307").Value = '3: 3607'
$ws.Range("H307").Value = 0
$ws.Range("I307").Value = '_x0003_32'
$ws.Range("J307").Value = 3
$ws.Range("K307").Value = 0.01567316232171778
$ws.Range("L307").Value = 'chen'
$ws.Range("M307").Value = '1/31/19 10:23:49'

# Row 308
$ws.Range("A308").Value = $bullet
$ws.Range("B308").Value = ''
$ws.Range("C308").Value = ''
$ws.Range("D308").Value = '4519'
$ws.Range("E308").Value = 'MIC'
$ws.Range("F308").Value = '3: 3615'
$ws.Range("G308").Value = '3: 3617'
$ws.Range("H308").Value = 0
$ws.Range("I308").Value = '_x0003_64'
$ws.Range("J308").Value = 3
$ws.Range("K308").Value = 0.01567316232171778
$ws.Range("L308").Value = 'chen'
$ws.Range("M308").Value = '1/31/19 10:23:56'

# Row 309
$ws.Range("A309").Value = $bullet
$ws.Range("B309").Value = ''
$ws.Range("C309").Value = ''
$ws.Range("D309").Value = '4519'
$ws.Range("E309").Value = 'MIC'
$ws.Range("F309").Value = '3: 3620'
$ws.Range("G309").Value = '3: 3622'
$ws.Range("H309").Value = 0
$ws.Range("I309").Value = '_x0003_64'
$ws.Range("J309").Value = 3
$ws.Range("K309").Value = 0.01567316232171778
$ws.Range("L309").Value = 'chen'
$ws.Range("M309").Value = '1/31/19 10:24:04'

# Row 310
$ws.Range("A310").Value = $bullet
$ws.Range("B310").Value = ''
$ws.Range("C310").Value = ''
$ws.Range("D310").Value = '4519'
$ws.Range("E310").Value = 'MIC'
$ws.Range("F310").Value = '3: 3631'
$ws.Range("G310").Value = '3: 3634'
$ws.Range("H310").Value = 0
$ws.Range("I310").Value = '_x0003_128'
$ws.Range("J310").Value = 4
$ws.Range("K310").Value = 0.02089754976229037
$ws.Range("L310").Value = 'chen'
$ws.Range("M310").Value = '1/31/19 10:24:33'

# Row 311
$ws.Range("A311").Value = $bullet
$ws.Range("B311").Value = ''
$ws.Range("C311").Value = ''
$ws.Range("D311").Value = '4519'
$ws.Range("E311").Value = 'MIC'
$ws.Range("F311").Value = '3: 3639'
$ws.Range("G311").Value = '3: 3641'
$ws.Range("H311").Value = 0
$ws.Range("I311").Value = '_x0003_32'
$ws.Range("J311").Value = 3
$ws.Range("K311").Value = 0.01567316232171778
$ws.Range("L311").Value = 'chen'
$ws.Range("M311").Value = '1/31/19 10:24:41'

# Row 312
$ws.Range("A312").Value = $bullet
$ws.Range("B312").Value = ''
$ws.Range("C312").Value = ''
$ws.Range("D312").Value = '4519'
$ws.Range("E312").Value = 'MIC'
$ws.Range("F312").Value = '3: 3644'
$ws.Range("G312").Value = '3: 3646'
$ws.Range("H312").Value = 0
$ws.Range("I312").Value = '_x0003_32'
$ws.Range("J312").Value = 3
$ws.Range("K312").Value = 0.01567316232171778
$ws.Range("L312").Value = 'chen'
$ws.Range("M312").Value = '1/31/19 10:24:46'

